# Apply updated numeric results from the parallel CPU run (CUDA work-in-progress)
# to the summary statistics on Sheet1 and Sheet2.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = 1.289713969429163
$ws.Range("C2").Value = 1.040549717476848
$ws.Range("D2").Value = 0.2491642519523148
$ws.Range("E2").Value = 0.1968626860971182
$ws.Range("F2").Value = 0.8641083595175394
$ws.Range("G2").Value = 2.161941513401254
$ws.Range("H2").Value = -0.00001752153088167524
$ws.Range("B3").Value = 0.0001617203364274568
$ws.Range("C3").Value = 0.0001095081546554957
$ws.Range("D3").Value = 0.00005587731278983738
$ws.Range("E3").Value = 0.00001637100151948144
$ws.Range("F3").Value = 0.0001083526254063959
$ws.Range("G3").Value = 0.0003378014666888245
$ws.Range("H3").Value = 0.0001285032358530241
$ws.Range("C4").Value = 0.6756238030766768
$ws.Range("D4").Value = 0.3474293626494357
$ws.Range("E4").Value = 0.1021510693381799
$ws.Range("F4").Value = 0.669999999999999
$ws.Range("G4").Value = 2.07943393139117
$ws.Range("H4").Value = 0.7993801889009886
$ws.Range("C5").Value = 0.9884494485792923
$ws.Range("D5").Value = 0.9561650976675447
$ws.Range("E5").Value = -0.07549240067504825
$ws.Range("F5").Value = 0.9999999999999998
$ws.Range("G5").Value = 0.8412700963316098
$ws.Range("H5").Value = 0.9485925684281818
$ws.Range("B6").Value = 0.9495053716747701
$ws.Range("C6").Value = 0.9705347381981705
$ws.Range("D6").Value = 0.8982627044924651
$ws.Range("E6").Value = 0.8649238768848985
$ws.Range("F6").Value = 0.9495053716747688
$ws.Range("G6").Value = 0.9922147479707973
$ws.Range("H6").Value = 0.8937818712849989
$ws.Range("B7").Value = 0.0001253922401745724
$ws.Range("C7").Value = 0.0001052405553871066
$ws.Range("D7").Value = 0.0002242584882890124
$ws.Range("E7").Value = 0.00008315952163916979
$ws.Range("F7").Value = 0.0001253922401745723
$ws.Range("G7").Value = 0.0001562488713588043
$ws.Range("H7").Value = -11.20102049865351

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("B2").Value = 0.0000434972889381401
$ws.Range("B3").Value = 0.00004333066796271318
$ws.Range("B4").Value = 0.00004333870639554767
$ws.Range("B5").Value = 7.585690812714347

